# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Binance" conversion lines inside A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$texto = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.78 = 26971.19 pesos`n✅ 26971.19 pesos = 6.75 = 960.71 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $texto

# --- tasas: update the N10/O10 and N12/O12 rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 147.5
$wsTasas.Range("O10").Value = 3978.25
$wsTasas.Range("N12").Value = 3994.96
$wsTasas.Range("O12").Value = 142.3
